$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.504.24"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.912.89"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "'325.44"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").Value = "'0.4810"
$ws.Range("D8").Value = "'0.4057"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "'0.08131"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "'1.011"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "'23.50"
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("D12").Value = "1.913.94"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "'5.985"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").Value = "'7.116"
$ws.Range("D15").Value = "'90.17"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.009"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.06776"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "'17.64"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Value = "'1.007"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "29.532.71"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "'5.620"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "'11.75"
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("D24").Value = "'2.186"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").Value = "2.123.77"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'155.30"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "'6.420"
$ws.Range("E27").Value = "  +4.59%  "
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "'2.099"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").Value = "'119.74"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").Value = "'1.027"
$ws.Range("E31").Value = "  -4.36%  "
$ws.Range("D32").Value = "'0.09529"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").Value = "'5.496"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").Value = "'3.571"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("D35").Value = "'1.386"
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "'0.06106"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "'1.179"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.940"
$ws.Range("E40").Value = "  -4.39%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'10.68"
$ws.Range("E41").Value = "  +5.23%  "
$ws.Range("D42").Value = "'0.1852"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("D43").Value = "'2.417"
$ws.Range("E43").Value = "  -5.43%  "
$ws.Range("D44").Value = "'1.285"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").Value = "'0.07639"
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("D47").Value = "'0.5572"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "'1.937"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").Value = "'116.11"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D50").Value = "'72.33"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("E51").Value = "  +1.74%  "

# Reset style on cells where we used a leading apostrophe to force text,
# so no stray number-format / quote-prefix style sticks to the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
